$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PegaTestData")

# Handle hard coded -API: replace the hard-coded campaign test value with CampaignCD
$ws.Range("A2").Value = "CampaignCD"

# Update the active selection to reflect where the edit was made
$ws.Activate()
$ws.Range("C6").Select()
